$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column B ("voucherDate") -------------------------------
# This shifts the existing voucherType/accountCode1/accountCode2/department/
# function columns one place to the right (B->C, C->D, D->E, E->F, F->G)
# while preserving every existing cell's value/type/style untouched.
$ws.Columns.Item(2).Insert()

# --- New header cell for the inserted column ------------------------------
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "voucherDate"

# --- New row 4: a sample "voucherDate" data row ---------------------------
$ws.Range("A4").NumberFormat = "General"
$ws.Range("A4").Value = "voucherDate"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "30/06/2016"

$ws.Range("C4").NumberFormat = "General"
$ws.Range("C4").Value = "Expense"

$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Value = 2101001
$ws.Range("D4").NumberFormat = "@"

$ws.Range("E4").NumberFormat = "General"
$ws.Range("E4").Value = 3501003
$ws.Range("E4").NumberFormat = "@"

$ws.Range("F4").NumberFormat = "General"
$ws.Range("F4").Value = "PUBLIC HEALTH AND SANITATION"

$ws.Range("G4").NumberFormat = "General"
$ws.Range("G4").Value = "Public Health"

$ws.Rows.Item(4).RowHeight = 18.7

# --- Column widths (character-width COM property rounds to whole pixels;
#     subtract the 5/6 "chrome" offset so the stored xlsx width lands as
#     close as possible to the target) -------------------------------------
$pad = 5/6
$ws.Columns.Item(1).ColumnWidth = 11.8 - $pad
$ws.Columns.Item(2).ColumnWidth = 13.2296296296296 - $pad
$ws.Columns.Item(3).ColumnWidth = 13.2296296296296 - $pad
$ws.Columns.Item(4).ColumnWidth = 14.7962962962963 - $pad
$ws.Columns.Item(5).ColumnWidth = 16.462962962963 - $pad
$ws.Columns.Item(6).ColumnWidth = 30.0851851851852 - $pad
$ws.Columns.Item(7).ColumnWidth = 25.7740740740741 - $pad

# --- Selection, matching the target sheet view ----------------------------
$ws.Range("E13").Select()
